$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 62.666668
$ws.Range("I4").Value = 59.5
$ws.Range("K4").Value = 59.5
$ws.Range("M4").Value = 54.5

$ws.Range("H7").Value = 2000
$ws.Range("I7").Value = 2000
$ws.Range("K7").Value = 2000
$ws.Range("M7").Value = -1888

$ws.Range("H10").Value = 2000
$ws.Range("I10").Value = 2000
$ws.Range("K10").Value = 2000
$ws.Range("M10").Value = -1707

$ws.Range("H13").Value = 1000
$ws.Range("I13").Value = 1000
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 1000
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = -831
$ws.Range("N13").ClearContents()

$ws.Range("H14").Value = 2000
$ws.Range("I14").Value = 2000
$ws.Range("K14").Value = 2000
$ws.Range("M14").Value = -1809

$ws.Range("H20").Value = 4499
$ws.Range("I20").Value = 4499
$ws.Range("K20").Value = 4499
$ws.Range("M20").Value = -4269

$ws.Range("H34").Value = 977.125
$ws.Range("I34").Value = 363.4
$ws.Range("J34").Value = 2000
$ws.Range("K34").Value = 363.4
$ws.Range("L34").Value = 2000
$ws.Range("M34").Value = -160.4
$ws.Range("N34").Value = -2406

$ws.Range("H35").Value = 4499
$ws.Range("I35").Value = 4499
$ws.Range("K35").Value = 4499
$ws.Range("M35").Value = -4120

$ws.Range("H36").Value = 977.125
$ws.Range("I36").Value = 363.4
$ws.Range("J36").Value = 2000
$ws.Range("K36").Value = 363.4
$ws.Range("L36").Value = 2000
$ws.Range("M36").Value = 351.6
$ws.Range("N36").Value = -3430

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H12").Value = 774.25
$ws.Range("J12").Value = 965.6667
$ws.Range("L12").Value = 965.6667
$ws.Range("N12").Value = -1311.6667

$ws.Range("H38").Value = 499.5
$ws.Range("J38").Value = 499.5
$ws.Range("L38").Value = 499.5
$ws.Range("N38").Value = -1433.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H36").Value = 9024
$ws.Range("I36").Value = 9024
$ws.Range("K36").Value = 9024
$ws.Range("M36").Value = -8636

$ws.Range("H40").Value = 9024
$ws.Range("I40").Value = 9024
$ws.Range("K40").Value = 9024
$ws.Range("M40").Value = -8864

$ws.Range("H122").Value = 1297.625
$ws.Range("I122").Value = 1025.8572
$ws.Range("J122").Value = 3200
$ws.Range("K122").Value = 3077.5716
$ws.Range("L122").Value = 9600
$ws.Range("M122").Value = -627.5715999999998
$ws.Range("N122").Value = -14500

$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H29").Value = 151
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 151
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 453
$ws.Range("M29").ClearContents()
$ws.Range("N29").Value = -1007

$ws.Range("H34").Value = 675
$ws.Range("I34").Value = 350
$ws.Range("K34").Value = 1050
$ws.Range("M34").Value = -966

$ws.Range("H36").Value = 50
$ws.Range("I36").Value = 50
$ws.Range("K36").Value = 150
$ws.Range("M36").Value = 19

$ws.Range("H43").Value = 4000
$ws.Range("J43").Value = 4000
$ws.Range("L43").Value = 12000
$ws.Range("N43").Value = -12228

$ws.Range("H47").Value = 401
$ws.Range("I47").Value = 401.5
$ws.Range("J47").Value = 400
$ws.Range("K47").Value = 1204.5
$ws.Range("L47").Value = 1200
$ws.Range("M47").Value = -773.5
$ws.Range("N47").Value = -2062

$ws.Range("H50").Value = 8
$ws.Range("I50").Value = 8
$ws.Range("K50").Value = 24
$ws.Range("M50").Value = 457

$ws.Range("H53").Value = 8
$ws.Range("I53").Value = 8
$ws.Range("K53").Value = 24
$ws.Range("M53").Value = 457

$ws.Range("H59").Value = 129.66667
$ws.Range("I59").Value = 129.66667
$ws.Range("K59").Value = 389.00001
$ws.Range("M59").Value = 150.99999

$ws.Range("H60").Value = 150
$ws.Range("I60").Value = 150
$ws.Range("K60").Value = 450
$ws.Range("M60").Value = -199

$ws.Range("H61").Value = 600
$ws.Range("I61").Value = 200
$ws.Range("J61").Value = 1000
$ws.Range("K61").Value = 600
$ws.Range("L61").Value = 3000
$ws.Range("M61").Value = -385
$ws.Range("N61").Value = -3430

$ws.Range("H121").Value = 500
$ws.Range("I121").Value = 500
$ws.Range("K121").Value = 1500
$ws.Range("M121").Value = -190

$ws.Range("H137").Value = 4000
$ws.Range("I137").Value = 3000
$ws.Range("K137").Value = 9000
$ws.Range("M137").Value = -3900

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 40000
$ws.Range("J26").Value = 40000
$ws.Range("L26").Value = 40000
$ws.Range("N26").Value = -40560

$ws.Range("H50").Value = 40000
$ws.Range("J50").Value = 40000
$ws.Range("L50").Value = 40000
$ws.Range("N50").Value = -40996

$ws.Range("H97").Value = 733.3333
$ws.Range("I97").Value = 733.3333
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 733.3333
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -237.3333
$ws.Range("N97").ClearContents()

$ws.Range("H132").Value = 9059.23
$ws.Range("I132").Value = 6706.364
$ws.Range("K132").Value = 20119.092
$ws.Range("M132").Value = -17589.092

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 26875
$ws.Range("J2").Value = 100000
$ws.Range("L2").Value = 100000
$ws.Range("N2").Value = -100224

$ws.Range("H9").Value = 415
$ws.Range("I9").Value = 282
$ws.Range("K9").Value = 282
$ws.Range("M9").Value = -58

$ws.Range("H12").Value = 1122
$ws.Range("I12").Value = 499
$ws.Range("K12").Value = 499
$ws.Range("M12").Value = -329

$ws.Range("H17").Value = 4009
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").ClearContents()

$ws.Range("H31").Value = 1700
$ws.Range("J31").Value = 1700
$ws.Range("L31").Value = 1700
$ws.Range("N31").Value = -2196

$ws.Range("H47").Value = 5000
$ws.Range("J47").Value = 5000
$ws.Range("L47").Value = 5000
$ws.Range("N47").Value = -5980

$ws.Range("H52").Value = 5000
$ws.Range("J52").Value = 5000
$ws.Range("L52").Value = 5000
$ws.Range("N52").Value = -5466

$ws.Range("H58").Value = 5834.3335
$ws.Range("I58").Value = 150
$ws.Range("J58").Value = 8676.5
$ws.Range("K58").Value = 150
$ws.Range("L58").Value = 8676.5
$ws.Range("M58").Value = 110
$ws.Range("N58").Value = -9196.5

$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("M136").ClearContents()
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H23").Value = 1474.5
$ws.Range("I23").Value = 1474.5
$ws.Range("K23").Value = 1474.5
$ws.Range("M23").Value = -1245.5

$ws.Range("H136").Value = 9355.799999999999
$ws.Range("I136").Value = 9355.799999999999
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 28067.4
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -25517.4
$ws.Range("N136").ClearContents()
